$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as literal text even when the string looks numeric,
# by temporarily switching the cell to Text format, writing the value, then
# restoring the original style so no stray formatting change is left behind.
function Set-TextValue($cell, $value) {
    $orig = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $orig
}

$ws.Range("D2").Value = '30.558.39'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '1.882.17'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  -0.04%  '
Set-TextValue $ws.Range("D5") '246.45'
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("E6").Value = '  -0.02%  '
Set-TextValue $ws.Range("D7") '0.4726'
$ws.Range("E7").Value = '  -0.12%  '
Set-TextValue $ws.Range("D8") '0.2884'
$ws.Range("E8").Value = '  -1.20%  '
Set-TextValue $ws.Range("D9") '0.06532'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("E10").Value = '  +0.48%  '
$ws.Range("B11").Value = 'Litecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D11") '100.70'
$ws.Range("E11").Value = '  +4.12%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D12") '0.7691'
$ws.Range("E12").Value = '  +4.59%  '
Set-TextValue $ws.Range("D13") '0.07824'
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D14").Value = '1.879.74'
$ws.Range("E14").Value = '  -0.68%  '
Set-TextValue $ws.Range("D15") '5.243'
$ws.Range("E15").Value = '  -0.10%  '
Set-TextValue $ws.Range("D16") '284.94'
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("D17").Value = '30.542.76'
Set-TextValue $ws.Range("D18") '13.18'
$ws.Range("E18").Value = '  -0.40%  '
Set-TextValue $ws.Range("D19") '0.000007519'
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '2.124.40'
$ws.Range("E21").Value = '  -0.74%  '
Set-TextValue $ws.Range("D22") '5.358'
$ws.Range("E22").Value = '  +0.72%  '
Set-TextValue $ws.Range("D23") '1.0000'
$ws.Range("E23").Value = '  -0.03%  '
Set-TextValue $ws.Range("D24") '6.395'
$ws.Range("E24").Value = '  +2.34%  '
Set-TextValue $ws.Range("D25") '9.131'
$ws.Range("E25").Value = '  -1.03%  '
Set-TextValue $ws.Range("D26") '162.63'
$ws.Range("E26").Value = '  -1.08%  '
Set-TextValue $ws.Range("D27") '19.07'
$ws.Range("E27").Value = '  +0.65%  '
Set-TextValue $ws.Range("D28") '1.913'
$ws.Range("E28").Value = '  -0.44%  '
Set-TextValue $ws.Range("D29") '0.09698'
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("E30").Value = '  -1.20%  '
$ws.Range("E31").Value = '  +0.87%  '
Set-TextValue $ws.Range("D32") '4.259'
$ws.Range("E32").Value = '  -1.00%  '
Set-TextValue $ws.Range("D33") '4.195'
$ws.Range("E33").Value = '  -0.07%  '
Set-TextValue $ws.Range("D34") '0.04837'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("E35").Value = '  +0.18%  '
Set-TextValue $ws.Range("D36") '0.6957'
$ws.Range("E36").Value = '  -0.33%  '
Set-TextValue $ws.Range("D37") '2.757'
$ws.Range("E37").Value = '  +1.21%  '
$ws.Range("E38").Value = '  +0.70%  '
Set-TextValue $ws.Range("D39") '2.882'
$ws.Range("E39").Value = '  +2.64%  '
Set-TextValue $ws.Range("D40") '76.39'
$ws.Range("E40").Value = '  +0.51%  '
Set-TextValue $ws.Range("D41") '6.275'
$ws.Range("E41").Value = '  -1.47%  '
Set-TextValue $ws.Range("D42") '1.977'
$ws.Range("E42").Value = '  -1.33%  '
Set-TextValue $ws.Range("D43") '0.4256'
$ws.Range("E43").Value = '  +0.09%  '
$ws.Range("E44").Value = '  -0.08%  '
Set-TextValue $ws.Range("D45") '0.8291'
$ws.Range("E45").Value = '  -0.75%  '
Set-TextValue $ws.Range("D46") '101.62'
$ws.Range("E46").Value = '  -0.08%  '
Set-TextValue $ws.Range("D47") '9.811'
$ws.Range("E47").Value = '  +2.87%  '
Set-TextValue $ws.Range("D48") '7.032'
$ws.Range("E48").Value = '  +0.03%  '
Set-TextValue $ws.Range("D49") '35.13'
$ws.Range("E49").Value = '  -1.58%  '
Set-TextValue $ws.Range("D50") '892.41'
$ws.Range("E50").Value = '  -2.83%  '
Set-TextValue $ws.Range("D51") '0.05760'
$ws.Range("E51").Value = '  +0.14%  '
